$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- helpers -----------------------------------------------------------

# Plain text (non-date-looking) values are stored as text automatically;
# just assign them.
function Set-PlainText($cell, [string]$value) {
    $cell.Value = $value
}

# Values that Excel's auto-detection would otherwise coerce (e.g. strings
# that look like dates, or an empty value) need to be forced to stay as
# literal text, then have the residual NumberFormat stripped back off so
# no spurious style gets left behind on the cell.
function Set-ForcedText($cell, [string]$value) {
    if ($value -eq "") {
        # A lone leading apostrophe forces an empty *text* entry (as
        # opposed to a genuinely blank cell).
        $cell.Value = "'"
    } else {
        $cell.NumberFormat = "@"
        $cell.Value = $value
    }
    $cell.ClearFormats()
}

# Column A reuses the bordered/centred numeric "id" style already present
# on the sheet (e.g. A2); copy/paste-format it onto the new id cell. The
# clipboard gets invalidated by other operations (ClearFormats etc.), so
# re-copy right before every paste.
function Set-IdCell($cell, $value) {
    $cell.Value = $value
    $ws.Cells.Item(2, 1).Copy()
    $cell.PasteSpecial(-4122)
}

$rows = @(
    @{
        row = 6
        A = 4
        B = "Q7747"
        C = "Vladimir"
        D = "Putin"
        E = "1952-10-07"
        F = ""
        G = "Saint Petersburg"
        H = "https://ru.wikipedia.org/wiki/%D0%9F%D1%83%D1%82%D0%B8%D0%BD,_%D0%92%D0%BB%D0%B0%D0%B4%D0%B8%D0%BC%D0%B8%D1%80_%D0%92%D0%BB%D0%B0%D0%B4%D0%B8%D0%BC%D0%B8%D1%80%D0%BE%D0%B2%D0%B8%D1%87"
        I = "2000-05-07"
        J = "2008-05-07"
        K = "Владимир Владимирович Путин российский государственный и политический деятель. Действующий президент Российской Федерации, председатель Государственного Совета Российской Федерации и Совета Безопасности Российской Федерации Верховный главнокомандующий Вооружёнными силами Российской Федерации с 7 мая 2012 года. Ранее занимал должность президента с 7 мая 2000 по 7 мая 2008 года, также в 19992000 и 20082012 годах занимал должность председателя правительства Российской Федерации. Фактически руководит Россией, согласно разным оценкам, с 1999 или с 2000 года. В сентябре 2017 года Путин стал самым долго правящим российским лидером со времён Иосифа Сталина.`n"
    },
    @{
        row = 7
        A = 5
        B = "Q23530"
        C = "Dmitry"
        D = "Medvedev"
        E = "1965-09-14"
        F = ""
        G = "Saint Petersburg"
        H = "https://ru.wikipedia.org/wiki/%D0%9C%D0%B5%D0%B4%D0%B2%D0%B5%D0%B4%D0%B5%D0%B2,_%D0%94%D0%BC%D0%B8%D1%82%D1%80%D0%B8%D0%B9_%D0%90%D0%BD%D0%B0%D1%82%D0%BE%D0%BB%D1%8C%D0%B5%D0%B2%D0%B8%D1%87"
        I = "2008-05-07"
        J = "2012-05-07"
        K = "иностранные`n"
    },
    @{
        row = 8
        A = 6
        B = "Q34453"
        C = "Boris"
        D = "Yeltsin"
        E = "1931-02-01"
        F = "2007-04-23"
        G = "Butka"
        H = "https://ru.wikipedia.org/wiki/%D0%95%D0%BB%D1%8C%D1%86%D0%B8%D0%BD,_%D0%91%D0%BE%D1%80%D0%B8%D1%81_%D0%9D%D0%B8%D0%BA%D0%BE%D0%BB%D0%B0%D0%B5%D0%B2%D0%B8%D1%87"
        I = "1991-07-10"
        J = "1999-12-31"
        K = "Борис Николаевич Ельцин, Бутка, Буткинский район, Уральская область, СССР 23 апреля 2007, Москва, Россия советский и российский партийный, государственный и политический деятель, первый всенародно избранный Президент Российской Федерации в ноябре 1991 июне 1992 года одновременно возглавлял правительство. С марта по май 1992 года исполнял обязанности министра обороны Российской Федерации.`n"
    }
)

foreach ($r in $rows) {
    $rowNum = $r.row

    Set-IdCell $ws.Cells.Item($rowNum, 1) $r.A

    Set-PlainText  $ws.Cells.Item($rowNum, 2) $r.B
    Set-PlainText  $ws.Cells.Item($rowNum, 3) $r.C
    Set-PlainText  $ws.Cells.Item($rowNum, 4) $r.D
    Set-ForcedText $ws.Cells.Item($rowNum, 5) $r.E
    Set-ForcedText $ws.Cells.Item($rowNum, 6) $r.F
    Set-PlainText  $ws.Cells.Item($rowNum, 7) $r.G
    Set-PlainText  $ws.Cells.Item($rowNum, 8) $r.H
    Set-ForcedText $ws.Cells.Item($rowNum, 9) $r.I
    Set-ForcedText $ws.Cells.Item($rowNum, 10) $r.J
    Set-PlainText  $ws.Cells.Item($rowNum, 11) $r.K

    # The biography text contains an embedded newline, which makes the
    # COM value-setter wrap the row to a custom height; put the row back
    # to its natural (non-custom) height.
    $ws.Rows.Item($rowNum).AutoFit()
}
